$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$cell = $ws.Range("O9")
$cell.Value = "P"
$cell.ReadingOrder = 1
Write-Host "done"
